# Separate out the lat and lon in data file
#
# The worksheet has a column F "lat_lon" with strings like "(6.342472437999014, 100.18982195228926)".
# We split this into two new numeric columns G ("lat") and H ("lon").

function TruncSig($s, $sig) {
    # Truncate (not round) the decimal string $s to $sig significant digits,
    # then parse it as a double. This mirrors the precision-loss behaviour
    # observed in the source data (values stored with 15 significant digits).
    $neg = $false
    if ($s.StartsWith("-")) {
        $neg = $true
        $s = $s.Substring(1)
    }
    $intPart = $s
    $fracPart = ""
    if ($s.Contains(".")) {
        $parts = $s.Split(".")
        $intPart = $parts[0]
        $fracPart = $parts[1]
    }
    $stripped = $intPart.TrimStart("0")
    if ($stripped -eq "") {
        # Value is < 1 in magnitude -- count significant digits starting
        # from the first nonzero fractional digit.
        $i = 0
        while ($i -lt $fracPart.Length -and $fracPart.Substring($i, 1) -eq "0") {
            $i = $i + 1
        }
        $totalFracLen = $i + $sig
        if ($totalFracLen -gt $fracPart.Length) {
            $totalFracLen = $fracPart.Length
        }
        $newFrac = $fracPart.Substring(0, $totalFracLen)
        $result = "0.$newFrac"
    } else {
        $intSig = $stripped.Length
        $remaining = $sig - $intSig
        if ($remaining -lt 0) {
            $remaining = 0
        }
        if ($remaining -gt $fracPart.Length) {
            $remaining = $fracPart.Length
        }
        $newFrac = $fracPart.Substring(0, $remaining)
        $result = "$intPart.$newFrac"
    }
    if ($neg) {
        $result = "-$result"
    }
    return [double]$result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns.
$ws.Range("G1").Value = "lat"
$ws.Range("H1").Value = "lon"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Fill in lat / lon for each data row by parsing the "(lat, lon)" text in column F.
for ($r = 2; $r -le 30; $r++) {
    $txt = $ws.Range("F$r").Text
    if ($txt.StartsWith("(")) {
        $inner = $txt.Trim("()")
        $parts = $inner.Split(",")
        $latStr = $parts[0].Trim()
        $lonStr = $parts[1].Trim()
        $lat = TruncSig $latStr 15
        $lon = TruncSig $lonStr 15

        $ws.Range("G$r").Value = $lat
        $ws.Range("H$r").Value = $lon

        $ws.Range("F$r").Copy()
        $ws.Range("G${r}:H$r").PasteSpecial(-4122)
    }
}

# Adjust column widths to roughly match the post-edit layout.
$ws.Columns.Item(1).ColumnWidth = 5.166666666666667
$ws.Columns.Item(2).ColumnWidth = 7.385416666666667
$ws.Columns.Item(3).ColumnWidth = 13.944010416666666
$ws.Columns.Item(6).ColumnWidth = 37.830729166666664
$ws.Columns.Item(7).ColumnWidth = 11.166666666666666
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666

$ws.Range("G20").Select()
